$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rng = $ws.Range("C2:C468")
$rng.Value2 = 45179
